$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 258; this pushes the previous rows 258-280 down to 259-281
# and copies formatting (e.g. the date number-format on column D) from the row above.
$ws.Rows.Item(258).Insert()

# Populate the new row 258 with the new weekly record.
$ws.Range("A258").Value = 3
$ws.Range("B258").Value = "Femacal de La Calera"
$ws.Range("C258").Value = "Coquimbo"
$ws.Range("D258").Value = 44578
$ws.Range("E258").Value = 5
$ws.Range("F258").Value = 100112040
$ws.Range("G258").Value = "Cilantro"
$ws.Range("H258").Value = "Sin especificar"
$ws.Range("I258").Value = "Primera"
$ws.Range("J258").Value = 130
$ws.Range("K258").Value = 4000
$ws.Range("L258").Value = 4500
$ws.Range("M258").Value = 4269
$ws.Range("N258").Value = "$/docena de atados (3 kilos)"
$ws.Range("O258").Value = "Provincia de Quillota"
$ws.Range("P258").Value = 1423
$ws.Range("Q258").Value = 3
$ws.Range("R258").Value = "Hortaliza"
